# Registration sheet ("Basic details" last-record readback) gets refreshed
# with a newly generated random test record (account data for a new test run).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# Firstname / Middlename / Lastname / Username / DateOfBirth / AlternateEmail
# are plain text already -> a normal .Value assignment keeps them as shared
# strings with their existing cell styles untouched.
$ws.Range("B2").Value = "fnameJZPB"
$ws.Range("B3").Value = "mnameygHA"
$ws.Range("B4").Value = "lnameDoCa"
$ws.Range("B5").Value = "test44332@yopmail.com"
$ws.Range("B6").Value = "25/01/1966"
$ws.Range("B8").Value = "testAlt6286@yopmail.com"

# AlternateNumber is an all-digit string ("5987784429"); a bare .Value
# assignment would be auto-coerced into a numeric cell. Force it to stay
# text (leading apostrophe), then restore the plain General-format style
# (copied from a sibling cell) so only the value itself changes.
$ws.Range("B7").Value = "'5987784429"
$ws.Range("B3").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
